$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge the new header cell range first (before copying formats) so merged-range auto border splitting doesn't occur
$mergeResult = $ws.Range("K1:M1").Merge()

# Copy formats from the first block (columns B:D) onto the new block (columns K:M) for all 8 rows;
# column B:D carries the correct per-row border pattern (header / sub-header / data rows incl. bottom row)
$ws.Range("B1:D8").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Content
$ws.Range("K1").Value = "Bosch 0280155968 420cc"
$ws.Range("K2").Value = "ON"
$ws.Range("L2").Value = "OFF"
$ws.Range("M2").Value = "LAG"

# Formulas for M column (K - L) rows 3-8 (K and L themselves stay blank - draft placeholders).
# Mirror the existing D/G/J columns exactly: row 3 holds an independent (non-shared) formula,
# while rows 4-8 form the shared-formula group, to match the workbook's established pattern.
$ws.Range("M3").Formula = "=K3-L3"
$ws.Range("M4:M8").Formula = "=K4-L4"

# Fix I8 border to match the rest of the bottom row pattern (bug fix in original commit):
# it should use the same "middle column, bottom row" style as C8/F8 (bottom border only)
$ws.Range("C8").Copy()
$ws.Range("I8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move active cell/selection similar to target
$selResult = $ws.Range("M10").Select()

Write-Output "done"
